$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: new LED driver connector line item -----------------------------
# Write the new shared-string values in the same order the source workbook
# picked them up in (datasheet URL, part number, description) so the
# sharedStrings table gets appended in that order.
$ws.Range("D12").Value = "http://www.on-shore.com/wp-content/uploads/2015/09/ostvnxxa150.pdf"
$ws.Range("C12").Value = "OSTVN05A150"
$ws.Range("B12").Value = "CONN TERM BLOCK 2.54MM 5POS PCB"
$ws.Range("E12").Value = 2.34
$ws.Range("F12").Value = 1
$ws.Range("G12").Formula = "=F12*E12"

$ws.Hyperlinks.Add($ws.Range("D12"), "http://www.on-shore.com/wp-content/uploads/2015/09/ostvnxxa150.pdf")

# Adding the hyperlink swaps in a fresh ad-hoc "Hyperlink" flavoured style;
# restore the plain bordered hyperlink-column look the rest of the table uses.
$ws.Range("D11").Copy()
$ws.Range("D12").PasteSpecial(-4122)

# --- Extend the blank rows down to row 15, keeping the same look ------------
# Row 13 used to be the last (blank) row, and carried the shared-formula
# anchor in G. Clone its formatting into the two new blank rows 14 & 15.
$ws.Range("A13:G13").Copy()
$ws.Range("A14:G15").PasteSpecial(-4122)

# G13 itself becomes a plain blank cell; the running total formula now lives
# in the new trailing rows instead.
$ws.Range("G13").ClearContents()
$ws.Range("G14").Formula = "=F14*E14"
$ws.Range("G15").Formula = "=F15*E15"

$ws.Range("D9").Select()
